# Add yecc figures and updated report
#
# table_1 previously held a small "Scenario" lookup block (A1:E5) that is no
# longer needed; it is replaced with a new 1-row header describing the
# study/model/objective/climate-change figure mapping. table_1 also becomes
# the active sheet/tab (it previously was table_2).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("table_1")

# Wipe out the old "Scenario" block entirely (A1:E5) - this drops the
# now-unused shared strings (Human Development, growth, Demographics,
# Scenario, Scenario 1-4).
$ws1.Range("A1:E5").ClearContents()

# New header row for table_1.
$ws1.Range("A1").Value = "study"
$ws1.Range("B1").Value = "Model"
$ws1.Range("C1").Value = "Objective"
$ws1.Range("D1").Value = "Climate change"

# table_1 becomes the selected/active sheet (was table_2 before).
$ws1.Range("D2").Select()
$ws1.Activate()
